$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the 2020 column (J) of data ------------------------------------

# Header cell J3: copy format from I3 (year header style), set value 2020
$ws.Range("I3").Copy()
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("J3").Value = 2020

# J4: hazardous waste generated, thousand tons -> numeric, right format
$ws.Range("I4").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("J4").Value = 11545.7
$ws.Range("J4").NumberFormat = "0.0"

# J5: resident population -> stored as text "1 754,6" in the source file
$ws.Range("I5").Copy()
$ws.Range("J5").PasteSpecial(-4122)
$ws.Range("J5").NumberFormat = "0.0"
$ws.Range("J5").HorizontalAlignment = -4152
$ws.Range("J5").Value = "1 754,6"

# J6: computed ratio for 2020 -> numeric, matches format of I6
$ws.Range("I6").Copy()
$ws.Range("J6").PasteSpecial(-4122)
$ws.Range("J6").Value = 6636.8

# --- Remove the two now-unused trailing blank rows ----------------------
$ws.Range("A27:H28").EntireRow.Delete()

# --- Restore the selected cell as recorded in the saved workbook --------
$ws.Range("G22").Select()
